# error solve ifrs list
# Updates the financial figures for rows 2-6, removes the now-unused
# J/O (non-controlling interest) columns on rows 2-5, and clears out
# rows 7-9 entirely (keeping only the A/B/C label columns).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Row 2
# ---------------------------------------------------------------
$ws.Range("D2").Value = 1209
$ws.Range("E2").Value = 49
$ws.Range("F2").Value = 49
$ws.Range("G2").Value = 82
$ws.Range("H2").Value = 76
$ws.Range("I2").Value = 76
$ws.Range("J2").ClearContents()
$ws.Range("K2").Value = 1565
$ws.Range("L2").Value = 808
$ws.Range("M2").Value = 757
$ws.Range("N2").Value = 757
$ws.Range("O2").ClearContents()
$ws.Range("P2").Value = 50
$ws.Range("Q2").Value = 51
$ws.Range("R2").Value = 78
$ws.Range("S2").Value = -119
$ws.Range("T2").Value = 3
$ws.Range("U2").Value = 48
$ws.Range("V2").Value = 646
$ws.Range("W2").Value = 4.07
$ws.Range("X2").Value = 6.3
$ws.Range("Y2").Value = 10.59
$ws.Range("Z2").Value = 4.48
$ws.Range("AA2").Value = 106.63
$ws.Range("AB2").Value = 1419.25
$ws.Range("AC2").Value = 759
$ws.Range("AD2").Value = 4.32
$ws.Range("AE2").Value = 7545
$ws.Range("AF2").Value = 0.43
$ws.Range("AG2").Value = 0
$ws.Range("AH2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AJ2").Value = 10024319

# ---------------------------------------------------------------
# Row 3
# ---------------------------------------------------------------
$ws.Range("D3").Value = 1217
$ws.Range("E3").Value = 52
$ws.Range("F3").Value = 52
$ws.Range("G3").Value = 45
$ws.Range("H3").Value = 33
$ws.Range("I3").Value = 33
$ws.Range("J3").ClearContents()
$ws.Range("K3").Value = 1508
$ws.Range("L3").Value = 619
$ws.Range("M3").Value = 890
$ws.Range("N3").Value = 890
$ws.Range("O3").ClearContents()
$ws.Range("P3").Value = 62
$ws.Range("Q3").Value = 138
$ws.Range("R3").Value = -12
$ws.Range("S3").Value = -111
$ws.Range("T3").Value = 10
$ws.Range("U3").Value = 128
$ws.Range("V3").Value = 436
$ws.Range("W3").Value = 4.23
$ws.Range("X3").Value = 2.68
$ws.Range("Y3").Value = 3.96
$ws.Range("Z3").Value = 2.12
$ws.Range("AA3").Value = 69.54000000000001
$ws.Range("AB3").Value = 1337.16
$ws.Range("AC3").Value = 299
$ws.Range("AD3").Value = 16.7
$ws.Range("AE3").Value = 7138
$ws.Range("AF3").Value = 0.7
$ws.Range("AG3").Value = 150
$ws.Range("AH3").Value = 3
$ws.Range("AI3").Value = 57.32
$ws.Range("AJ3").Value = 12447744

# ---------------------------------------------------------------
# Row 4
# ---------------------------------------------------------------
$ws.Range("D4").Value = 1230
$ws.Range("E4").Value = 109
$ws.Range("F4").Value = 109
$ws.Range("G4").Value = 106
$ws.Range("H4").Value = 83
$ws.Range("I4").Value = 83
$ws.Range("J4").ClearContents()
$ws.Range("K4").Value = 1493
$ws.Range("L4").Value = 538
$ws.Range("M4").Value = 956
$ws.Range("N4").Value = 956
$ws.Range("O4").ClearContents()
$ws.Range("P4").Value = 62
$ws.Range("Q4").Value = 109
$ws.Range("R4").Value = -1
$ws.Range("S4").Value = -99
$ws.Range("T4").Value = 3
$ws.Range("U4").Value = 106
$ws.Range("V4").Value = 361
$ws.Range("W4").Value = 8.890000000000001
$ws.Range("X4").Value = 6.76
$ws.Range("Y4").Value = 9.01
$ws.Range("Z4").Value = 5.54
$ws.Range("AA4").Value = 56.26
$ws.Range("AB4").Value = 1440.48
$ws.Range("AC4").Value = 667
$ws.Range("AD4").Value = 8.19
$ws.Range("AE4").Value = 7668
$ws.Range("AF4").Value = 0.71
$ws.Range("AG4").Value = 200
$ws.Range("AH4").Value = 3.66
$ws.Range("AI4").Value = 30
$ws.Range("AJ4").Value = 12447744

# ---------------------------------------------------------------
# Row 5
# ---------------------------------------------------------------
$ws.Range("D5").Value = 1283
$ws.Range("E5").Value = 115
$ws.Range("F5").Value = 115
$ws.Range("G5").Value = 129
$ws.Range("H5").Value = 95
$ws.Range("I5").Value = 95
$ws.Range("J5").ClearContents()
$ws.Range("K5").Value = 1584
$ws.Range("L5").Value = 559
$ws.Range("M5").Value = 1025
$ws.Range("N5").Value = 1025
$ws.Range("O5").ClearContents()
$ws.Range("P5").Value = 62
$ws.Range("Q5").Value = -3
$ws.Range("R5").Value = -21
$ws.Range("S5").Value = 24
$ws.Range("T5").Value = 4
$ws.Range("U5").Value = -7
$ws.Range("V5").Value = 407
$ws.Range("W5").Value = 8.93
$ws.Range("X5").Value = 7.41
$ws.Range("Y5").Value = 9.59
$ws.Range("Z5").Value = 6.17
$ws.Range("AA5").Value = 54.49
$ws.Range("AB5").Value = 1552.93
$ws.Range("AC5").Value = 762
$ws.Range("AD5").Value = 7.02
$ws.Range("AE5").Value = 8225
$ws.Range("AF5").Value = 0.65
$ws.Range("AG5").Value = 270
$ws.Range("AH5").Value = 5.05
$ws.Range("AI5").Value = 35.42
$ws.Range("AJ5").Value = 12447744

# ---------------------------------------------------------------
# Row 6 (J6/O6 were already absent, no ClearContents needed there)
# ---------------------------------------------------------------
$ws.Range("D6").Value = 1201
$ws.Range("E6").Value = 73
$ws.Range("F6").Value = 73
$ws.Range("G6").Value = 77
$ws.Range("H6").Value = 56
$ws.Range("I6").Value = 56
$ws.Range("K6").Value = 1607
$ws.Range("L6").Value = 560
$ws.Range("M6").Value = 1047
$ws.Range("N6").Value = 1047
$ws.Range("P6").Value = 62
$ws.Range("Q6").Value = 20
$ws.Range("R6").Value = 1
$ws.Range("S6").Value = -26
$ws.Range("T6").Value = 5
$ws.Range("U6").Value = 16
$ws.Range("V6").Value = 413
$ws.Range("W6").Value = 6.09
$ws.Range("X6").Value = 4.66
$ws.Range("Y6").Value = 5.4
$ws.Range("Z6").Value = 3.51
$ws.Range("AA6").Value = 53.47
$ws.Range("AB6").Value = 1588.74
$ws.Range("AC6").Value = 449
$ws.Range("AD6").Value = 10.86
$ws.Range("AE6").Value = 8402
$ws.Range("AF6").Value = 0.58
$ws.Range("AG6").Value = 150
$ws.Range("AH6").Value = 3.08
$ws.Range("AI6").Value = 33.4
$ws.Range("AJ6").Value = 12447744

# ---------------------------------------------------------------
# Rows 7-9: remove all projected-year figures, keep only A/B/C labels
# ---------------------------------------------------------------
$ws.Range("D7:AJ7").ClearContents()
$ws.Range("D8:AJ8").ClearContents()
$ws.Range("D9:AJ9").ClearContents()

Write-Output "done"
